$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.693.83"
$ws.Range("E2").Value = "  +1.12%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.645.28"
$ws.Range("E3").Value = "  +0.39%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
$ws.Range("D5").Value = "212.93"
$ws.Range("E5").Value = "  +0.85%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.528"
$ws.Range("E6").Value = "  -0.59%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.01%  "

# Row 8 - Solana
$ws.Range("D8").Value = "23.39"
$ws.Range("E8").Value = "  +1.15%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +1.19%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.66%  "

# Row 11 - TRON
$ws.Range("D11").Value = "0.0893"
$ws.Range("E11").Value = "  +0.37%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.877.84"
$ws.Range("E12").Value = "  +0.36%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.639.32"
$ws.Range("E13").Value = "  +0.02%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "4.05"
$ws.Range("E14").Value = "  +0.83%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.562"
$ws.Range("E15").Value = "  +1.21%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "64.75"
$ws.Range("E16").Value = "  +0.89%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "27.671.12"
$ws.Range("E17").Value = "  +1.13%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "231.46"
$ws.Range("E18").Value = "  +0.57%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0725"
$ws.Range("E19").Value = "  +0.97%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "7.66"
$ws.Range("E20").Value = "  +1.76%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "4.30"
$ws.Range("E22").Value = "  +0.05%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  +6.98%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -2.73%  "

# Row 25 - Monero
$ws.Range("D25").Value = "149.89"
$ws.Range("E25").Value = "  +1.24%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "6.96"
$ws.Range("E26").Value = "  +0.19%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -1.21%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "15.68"
$ws.Range("E28").Value = "  +1.11%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.01%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +1.08%  "

# Row 31 - Hedera
$ws.Range("D31").Value = "0.0488"
$ws.Range("E31").Value = "  +0.99%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "3.30"
$ws.Range("E32").Value = "  +0.88%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.447.13"
$ws.Range("E33").Value = "  +3.02%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +0.91%  "

# Row 35 - LidoDAOToken
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  +1.11%  "

# Row 36 - HuobiToken
$ws.Range("E36").Value = "  -1.00%  "

# Row 37 - ImmutableX
$ws.Range("D37").Value = "0.571"
$ws.Range("E37").Value = "  +1.73%  "

# Row 38 - ARBITRUM
$ws.Range("D38").Value = "0.884"
$ws.Range("E38").Value = "  +0.69%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +0.62%  "

# Row 40 - TrustWalletToken
$ws.Range("D40").Value = "0.887"
$ws.Range("E40").Value = "  +12.70%  "

# Row 41 - WEMIXToken
$ws.Range("E41").Value = "  +0.45%  "

# Row 42 - was PaxDollar, now Aave
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "69.47"
$ws.Range("E42").Value = "  +7.83%  "

# Row 43 - was FraxShare, now PaxDollar
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.04%  "

# Row 44 - was Aave, now FraxShare
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.63"
$ws.Range("E44").Value = "  +3.28%  "

# Row 45 - mCoin
$ws.Range("E45").Value = "  +0.38%  "

# Row 46 - MXToken
$ws.Range("E46").Value = "  +0.47%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.787.35"
$ws.Range("E47").Value = "  +0.36%  "

# Row 48 - RenderToken
$ws.Range("D48").Value = "1.73"
$ws.Range("E48").Value = "  +5.89%  "

# Row 49 - Quant
$ws.Range("D49").Value = "85.78"
$ws.Range("E49").Value = "  -1.69%  "

# Row 50 - Algorand
$ws.Range("D50").Value = "0.0990"
$ws.Range("E50").Value = "  +0.24%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "7.78"
$ws.Range("E51").Value = "  +1.51%  "
